# Planning-Assignment.xlsx edit: "Remove error folder Admin-site"
#
# 1. Shift the "feature name" column (A:C) content for rows 18-23 down by one row,
#    and put a new "Sắp xếp sản phẩm" entry at row 18 (columns D:G keep their
#    original row positions - this was a partial-column insert, not a full row
#    insert).
# 2. Mark most feature rows as "Complete" in column C.
# 3. Drop the old mini Todolist entries (Code product-detail / check / Xử lý file
#    database) that lived under the "Todolist - (10-10-2022)" header, leaving just
#    the Todo/Done headers.
# 4. Move the active-cell selection to C26 (where the sheet now ends).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift column A:C content down by one row, starting at row 23 down to 18 ---
for ($r = 23; $r -ge 18; $r--) {
    $ws.Range("A" + ($r + 1)).Value = $ws.Range("A" + $r).Value2
    $ws.Range("B" + ($r + 1)).Value = $ws.Range("B" + $r).Value2
}
$ws.Range("A18").Value = $null

# --- 2. Mark rows as Complete in column C ---
$completeRows = @(3,5,6,7,9,10,11,12,13,14,15,16,17,18,19,20,22,23)
foreach ($r in $completeRows) {
    $ws.Range("C" + $r).Value = "Complete"
}

$ws.Range("B18").Value = "Sắp xếp sản phẩm "

# --- 3. Remove the old Todolist sub-items (rows 28 and 29) ---
$ws.Rows.Item(29).Delete()
$ws.Rows.Item(28).Delete()

# --- 4. Update the selection shown in the saved view ---
$ws.Range("C26").Select()
